# Fruta / hortaliza, semanal
# Insert a new weekly record above the current first row of the
# "Vega Monumental Concepción" / Zapallo italiano block (row 199),
# pushing the existing rows 199-206 down to 200-207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 199:206 down by inserting a new blank row at 199.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with this week's data.
$ws.Range("A199").Value = 11
$ws.Range("B199").Value = "Vega Monumental Concepción"
$ws.Range("C199").Value = "Bíobío"
$ws.Range("D199").Value = 45008
$ws.Range("E199").Value = 8
$ws.Range("F199").Value = 100112032
$ws.Range("G199").Value = "Zapallo italiano"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 100
$ws.Range("K199").Value = 8000
$ws.Range("L199").Value = 8500
$ws.Range("M199").Value = 8250
$ws.Range("N199").Value = "$/caja 50 unidades"
$ws.Range("O199").Value = "Región Metropolitana"
$ws.Range("P199").Value = 165
$ws.Range("Q199").Value = 50
$ws.Range("R199").Value = "Hortaliza"
